{"js": "const pairs = [\n  [\"Lomnxeba wefowuni uzakusithatha nje malunga nemizuzu elishumi elinesihlanu. Emva kokuba sigqibile, uzakufumana i-15 yeerandi zedata yamakhasi onxibelelwano yokuthi enkosi. Singathanda ukubuza imibuzo embalwa malunga nobom bakho kunye nobudlelwane bakho no X. Ndisebenza kwi projekthi enikeza ngenkxaso kubazali kunye nabanonopheli ngo WhatsApp. We are looking at how it's working and how the children have found it. Sifuna ukuva ukuba luthini uluvo lwakho ngayo. Akukho zimpendulo zilungileyo okanye zingalunganga apha. Nyaniseka qha, sifuna ukuva yintoni \", \"Lomnxeba wefowuni uzakusithatha nje malunga nemizuzu elishumi elinesihlanu. Emva kokuba sigqibile, uzakufumana i-15 yeerandi zedata yamakhasi onxibelelwano yokuthi enkosi. Singathanda ukubuza imibuzo embalwa malunga nobom bakho kunye nobudlelwane bakho no X. Ndisebenza kwi projekthi enikeza ngenkxaso kubazali kunye nabanonopheli ngo WhatsApp. Sijonga indlela esebenza ngayo kwaye indlela abantwana abayifumanisa ingayo. Sifuna ukuva ukuba luthini uluvo lwakho ngayo. Akukho zimpendulo zilungileyo okanye zingalunganga apha. Nyaniseka qha, sifuna ukuva yintoni \"],\n  [\"Ndifuna ukwabelana ngeempendulo ozinikezayo, kodwa xa ndisabelana ngazo, ndizakukunika igama elahlukileyo. Sifuna ukugcina igama lakho kunye nesazisi siyimfihlo ukuze sisebenzise elinye igama xa sisabelana nantoni na. Would you like to choose the name we use for you?   \", \"Ndifuna ukwabelana ngeempendulo ozinikezayo, kodwa xa ndisabelana ngazo, ndizakukunika igama elahlukileyo. Sifuna ukugcina igama lakho kunye nesazisi siyimfihlo ukuze sisebenzise elinye igama xa sisabelana nantoni na. Ingaba ungathanda ukukhetha igama esinolu sebenzisa for wena?   \"],\n  [\"Ukuba kulungile kuwe ndizakuwushicilela lomnxeba walefowuni, ukuze ndikhumbule iimpendulo zakho kamva. Only the people working with me on this project will hear what you say, but we won\u2019t share your real name with anyone. Do I have permission to record? \", \"Ukuba kulungile kuwe ndizakuwushicilela lomnxeba walefowuni, ukuze ndikhumbule iimpendulo zakho kamva. Ngabantu abasebenza nam kule projekthi kuphela abazakuva into oyitshoyo, kodwa asizokwabelana ngegama lakho lokwenene nakubani na. Ingaba ndinayo imvume yokurekhoda? \"],\n  [\"Ndizakubuza imibuzo embalwa. Lena yonke yimibuzo yokhetho oluninzi. Oku kuthetha ukuba xa uphendula, uzakukhetha kuphela impendulo engcono. Iinketho zezi 'Nakanye', 'Ngamanye amaxesha' okanye 'Rhoqo'. Yonke imibuzo ibibuza malunga nokuba kukangaphi kwiveki ephelileyo into isenzeka. You will then choose the best answer from \u2018Never\u2019, \u2018Sometimes\u2019 or \u2018Often\u2019. Sifuna kuphela ukuba uphendule malunga ne neveki ephelileyo okanye intsuku ezisixhenxe. Namhlanje yi _____, oko kuthetha ukuba ndibuza kungangaphi into isenzeka usukela okokugqibela ____. Then you can answer \u2018Never\u2019 if it hasn\u2019t happened at all since last ____, \u2018Often\u2019 if it happened almost every day, or about 5 times in the last week, and \u2018Sometimes\u2019 if it's somewhere between, such as 1 or 2 times. Ingaba lento iyavakala? (linda impendulo) \", \"Ndizakubuza imibuzo embalwa. Lena yonke yimibuzo yokhetho oluninzi. Oku kuthetha ukuba xa uphendula, uzakukhetha kuphela impendulo engcono. Iinketho zezi 'Nakanye', 'Ngamanye amaxesha' okanye 'Rhoqo'. Yonke imibuzo ibibuza malunga nokuba kukangaphi kwiveki ephelileyo into isenzeka. Uya kuthi emva koko ukhethe eyona mpendulo ingcono ukusuka ku-'Zange', 'Ngamanye amaxesha' okanye 'Rhoqo'. Sifuna kuphela ukuba uphendule malunga ne neveki ephelileyo okanye intsuku ezisixhenxe. Namhlanje yi _____, oko kuthetha ukuba ndibuza kungangaphi into isenzeka usukela okokugqibela ____. Emva koko unokuphendula ngokuthi 'Zange' ukuba khange yenzeke kwaphela ukusukela okokugqibela ____, 'Rhoqo' ukuba yenzeka phantse yonke imihla, okanye malunga namaxesha amahlanu kwiveki ephelileyo, kwaye 'Ngamanye amaxesha' ukuba iphakathi, njengakanye okanye kabini amaxesha. Ingaba lento iyavakala? (linda impendulo) \"],\n  [\"Okay, are you ready? Ndizakuhamba umbuzo ngamnye nawe kancinci-kancinci. Ndicela undazise ukuba awuyiqondi into endiyibuzayo, okanye ukuba unayo nayiphi na imibuzo ngelixa sihambayo.\", \"Kulungile, ulungele? Ndizakuhamba umbuzo ngamnye nawe kancinci-kancinci. Ndicela undazise ukuba awuyiqondi into endiyibuzayo, okanye ukuba unayo nayiphi na imibuzo ngelixa sihambayo.\"],\n  [\"Since last _____, how often did your x talk to you about using your cell phone or being on the internet? Zange, Ngamanye amaxesha okanye Rhoqo.\", \"Ukugqibela kwakho _____, kungangaphi u-x wakho ethetha nawe malunga nokusebenzisa unomyayi wakho okanye ukuba kwi-intanethi? Zange, Ngamanye amaxesha okanye Rhoqo.\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match, found \" + results.items.length + \" for: \" + oldText.substring(0, 60));\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@('Lomnxeba wefowuni uzakusithatha nje malunga nemizuzu elishumi elinesihlanu. Emva kokuba sigqibile, uzakufumana i-15 yeerandi zedata yamakhasi onxibelelwano yokuthi enkosi. Singathanda ukubuza imibuzo embalwa malunga nobom bakho kunye nobudlelwane bakho no X. Ndisebenza kwi projekthi enikeza ngenkxaso kubazali kunye nabanonopheli ngo WhatsApp. We are looking at how it''s working and how the children have found it. Sifuna ukuva ukuba luthini uluvo lwakho ngayo. Akukho zimpendulo zilungileyo okanye zingalunganga apha. Nyaniseka qha, sifuna ukuva yintoni ', 'Lomnxeba wefowuni uzakusithatha nje malunga nemizuzu elishumi elinesihlanu. Emva kokuba sigqibile, uzakufumana i-15 yeerandi zedata yamakhasi onxibelelwano yokuthi enkosi. Singathanda ukubuza imibuzo embalwa malunga nobom bakho kunye nobudlelwane bakho no X. Ndisebenza kwi projekthi enikeza ngenkxaso kubazali kunye nabanonopheli ngo WhatsApp. Sijonga indlela esebenza ngayo kwaye indlela abantwana abayifumanisa ingayo. Sifuna ukuva ukuba luthini uluvo lwakho ngayo. Akukho zimpendulo zilungileyo okanye zingalunganga apha. Nyaniseka qha, sifuna ukuva yintoni ')\n  ,@('Ndifuna ukwabelana ngeempendulo ozinikezayo, kodwa xa ndisabelana ngazo, ndizakukunika igama elahlukileyo. Sifuna ukugcina igama lakho kunye nesazisi siyimfihlo ukuze sisebenzise elinye igama xa sisabelana nantoni na. Would you like to choose the name we use for you?   ', 'Ndifuna ukwabelana ngeempendulo ozinikezayo, kodwa xa ndisabelana ngazo, ndizakukunika igama elahlukileyo. Sifuna ukugcina igama lakho kunye nesazisi siyimfihlo ukuze sisebenzise elinye igama xa sisabelana nantoni na. Ingaba ungathanda ukukhetha igama esinolu sebenzisa for wena?   ')\n  ,@('Ukuba kulungile kuwe ndizakuwushicilela lomnxeba walefowuni, ukuze ndikhumbule iimpendulo zakho kamva. Only the people working with me on this project will hear what you say, but we won\u2019t share your real name with anyone. Do I have permission to record? ', 'Ukuba kulungile kuwe ndizakuwushicilela lomnxeba walefowuni, ukuze ndikhumbule iimpendulo zakho kamva. Ngabantu abasebenza nam kule projekthi kuphela abazakuva into oyitshoyo, kodwa asizokwabelana ngegama lakho lokwenene nakubani na. Ingaba ndinayo imvume yokurekhoda? ')\n  ,@('Ndizakubuza imibuzo embalwa. Lena yonke yimibuzo yokhetho oluninzi. Oku kuthetha ukuba xa uphendula, uzakukhetha kuphela impendulo engcono. Iinketho zezi ''Nakanye'', ''Ngamanye amaxesha'' okanye ''Rhoqo''. Yonke imibuzo ibibuza malunga nokuba kukangaphi kwiveki ephelileyo into isenzeka. You will then choose the best answer from \u2018Never\u2019, \u2018Sometimes\u2019 or \u2018Often\u2019. Sifuna kuphela ukuba uphendule malunga ne neveki ephelileyo okanye intsuku ezisixhenxe. Namhlanje yi _____, oko kuthetha ukuba ndibuza kungangaphi into isenzeka usukela okokugqibela ____. Then you can answer \u2018Never\u2019 if it hasn\u2019t happened at all since last ____, \u2018Often\u2019 if it happened almost every day, or about 5 times in the last week, and \u2018Sometimes\u2019 if it''s somewhere between, such as 1 or 2 times. Ingaba lento iyavakala? (linda impendulo) ', 'Ndizakubuza imibuzo embalwa. Lena yonke yimibuzo yokhetho oluninzi. Oku kuthetha ukuba xa uphendula, uzakukhetha kuphela impendulo engcono. Iinketho zezi ''Nakanye'', ''Ngamanye amaxesha'' okanye ''Rhoqo''. Yonke imibuzo ibibuza malunga nokuba kukangaphi kwiveki ephelileyo into isenzeka. Uya kuthi emva koko ukhethe eyona mpendulo ingcono ukusuka ku-''Zange'', ''Ngamanye amaxesha'' okanye ''Rhoqo''. Sifuna kuphela ukuba uphendule malunga ne neveki ephelileyo okanye intsuku ezisixhenxe. Namhlanje yi _____, oko kuthetha ukuba ndibuza kungangaphi into isenzeka usukela okokugqibela ____. Emva koko unokuphendula ngokuthi ''Zange'' ukuba khange yenzeke kwaphela ukusukela okokugqibela ____, ''Rhoqo'' ukuba yenzeka phantse yonke imihla, okanye malunga namaxesha amahlanu kwiveki ephelileyo, kwaye ''Ngamanye amaxesha'' ukuba iphakathi, njengakanye okanye kabini amaxesha. Ingaba lento iyavakala? (linda impendulo) ')\n  ,@('Okay, are you ready? Ndizakuhamba umbuzo ngamnye nawe kancinci-kancinci. Ndicela undazise ukuba awuyiqondi into endiyibuzayo, okanye ukuba unayo nayiphi na imibuzo ngelixa sihambayo.', 'Kulungile, ulungele? Ndizakuhamba umbuzo ngamnye nawe kancinci-kancinci. Ndicela undazise ukuba awuyiqondi into endiyibuzayo, okanye ukuba unayo nayiphi na imibuzo ngelixa sihambayo.')\n  ,@('Since last _____, how often did your x talk to you about using your cell phone or being on the internet? Zange, Ngamanye amaxesha okanye Rhoqo.', 'Ukugqibela kwakho _____, kungangaphi u-x wakho ethetha nawe malunga nokusebenzisa unomyayi wakho okanye ukuba kwi-intanethi? Zange, Ngamanye amaxesha okanye Rhoqo.')\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.Text = $oldText\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $found = $find.Execute()\n  if (-not $found) {\n    throw \"Find failed for: $($oldText.Substring(0, [Math]::Min(60, $oldText.Length)))\"\n  }\n  $find.Parent.Text = $newText\n}"}
